$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the last existing data row (24) down onto the
# three new rows (25-27) so the new cells reuse the same style records
# (border/fill/font) as the rest of the "perse" (珀斯古城) scene block.
$ws.Range("A24:T24").Copy()
$ws.Range("A25:T27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the Id / Type / Level / ReviveScene numbers first (these are plain
# numbers, so they do not touch the shared-string table and their order
# relative to the text below does not matter).
$ws.Range("A25").Value = 13020011
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = 13010007

$ws.Range("A26").Value = 13020012
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = 13010007

$ws.Range("A27").Value = 13020013
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 8
$ws.Range("E27").Value = 13010007

# Text/lookup cells - set in the same order the author typed them so the
# shared-string table grows with the same new-entry order as the source
# workbook: 古城大厅, persepalace1, persepalace2, persepalace3, 古城外围, 古城大殿.
$ws.Range("B26").Value = "古城大厅"
$ws.Range("H25").Value = "trees;4|manflower;2|portal;1|sandland;2"
$ws.Range("P25").Value = "persepalace1"
$ws.Range("Q25").Value = "persepalace1"
$ws.Range("P26").Value = "persepalace2"
$ws.Range("Q26").Value = "persepalace2"
$ws.Range("P27").Value = "persepalace3"
$ws.Range("Q27").Value = "persepalace3"
$ws.Range("B25").Value = "古城外围"
$ws.Range("B27").Value = "古城大殿"

# Grow the "表1" table (ListObject) and its autofilter to cover the new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A3:T27"))

# Extend the "is-zero" conditional-format highlight (I4:N24 -> I4:N27) so it
# keeps covering the QPortal..QAngel flag columns for the new rows too.
$fc = $ws.Range("I4:N24").FormatConditions
$fc.Item(1).ModifyAppliesToRange($ws.Range("I4:N27"))

# Move the active cell/selection the way the author left it.
$ws.Range("G16").Select()
